$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 490.83334
$ws.Range("I28").Value = 490.83334
$ws.Range("K28").Value = 490.83334
$ws.Range("M28").Value = -5.833340000000021
# Row 40
$ws.Range("H40").Value = 2116.111
$ws.Range("I40").Value = 2105.625
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 2105.625
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -1930.625
$ws.Range("N40").Value = -2550
# Row 62
$ws.Range("H62").Value = 9931.666999999999
$ws.Range("I62").Value = 9897.5
$ws.Range("K62").Value = 9897.5
$ws.Range("M62").Value = -9273.5
# Row 65
$ws.Range("H65").Value = 9931.666999999999
$ws.Range("I65").Value = 9897.5
$ws.Range("K65").Value = 49487.5
$ws.Range("M65").Value = -46367.5
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 111
$ws.Range("H111").Value = 2881.6667
$ws.Range("J111").Value = 3266.6667
$ws.Range("L111").Value = 9800.000100000001
$ws.Range("N111").Value = -15934.0001
# Row 137
$ws.Range("H137").Value = 2912.7778
$ws.Range("I137").Value = 2807.8333
$ws.Range("K137").Value = 8423.499899999999
$ws.Range("M137").Value = -5873.499899999999
# Row 138
$ws.Range("H138").Value = 10620.571
$ws.Range("J138").Value = 10620.571
$ws.Range("L138").Value = 31861.713
$ws.Range("N138").Value = -42141.713

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 24613.273
$ws.Range("I32").Value = 23074.7
$ws.Range("J32").Value = 39999
$ws.Range("K32").Value = 23074.7
$ws.Range("L32").Value = 39999
$ws.Range("M32").Value = -22787.7
$ws.Range("N32").Value = -40573
# Row 45
$ws.Range("H45").Value = 1862.25
$ws.Range("I45").Value = 2249.5
$ws.Range("J45").Value = 1475
$ws.Range("K45").Value = 2249.5
$ws.Range("L45").Value = 1475
$ws.Range("M45").Value = -1872.5
$ws.Range("N45").Value = -2229
# Row 122
$ws.Range("H122").Value = 250000000
$ws.Range("I122").Value = 250000000
$ws.Range("K122").Value = 750000000
$ws.Range("M122").Value = -749997550
# Row 132
$ws.Range("H132").Value = 6593.857
$ws.Range("I132").Value = 6518.077
$ws.Range("J132").Value = 7579
$ws.Range("K132").Value = 19554.231
$ws.Range("L132").Value = 22737
$ws.Range("M132").Value = -17024.231
$ws.Range("N132").Value = -27797

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 19236
$ws.Range("I82").Value = 8159.6
$ws.Range("K82").Value = 8159.6
$ws.Range("M82").Value = -7776.6
# Row 85
$ws.Range("H85").Value = 19236
$ws.Range("I85").Value = 8159.6
$ws.Range("K85").Value = 8159.6
$ws.Range("M85").Value = -6833.6
# Row 105
$ws.Range("H105").Value = 1094.8572
$ws.Range("I105").Value = 1010.6
$ws.Range("J105").Value = 1305.5
$ws.Range("K105").Value = 1010.6
$ws.Range("L105").Value = 1305.5
$ws.Range("M105").Value = 736.4
$ws.Range("N105").Value = -4799.5
# Row 107
$ws.Range("H107").Value = 1227.7142
$ws.Range("I107").Value = 1323.5
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1323.5
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 596.5
$ws.Range("N107").Value = -4940

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 2000
$ws.Range("J47").Value = 2000
$ws.Range("L47").Value = 2000
$ws.Range("N47").Value = -3132
# Row 107
$ws.Range("H107").Value = 324.33334
$ws.Range("I107").Value = 230
$ws.Range("K107").Value = 230
$ws.Range("M107").Value = 1690

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 66822.664
$ws.Range("J2").Value = 161.33333
$ws.Range("L2").Value = 967.9999799999999
$ws.Range("N2").Value = -1193.99998
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 70
$ws.Range("H70").Value = 3992.3333
$ws.Range("I70").Value = 2989
$ws.Range("K70").Value = 8967
$ws.Range("M70").Value = -8652
# Row 73
$ws.Range("H73").Value = 3992.3333
$ws.Range("I73").Value = 2989
$ws.Range("K73").Value = 8967
$ws.Range("M73").Value = -7875
# Row 80
$ws.Range("H80").Value = 1475
$ws.Range("I80").Value = 950
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -1914
$ws.Range("N80").Value = -7872
# Row 83
$ws.Range("H83").Value = 1475
$ws.Range("I83").Value = 950
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 8550
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -3870
$ws.Range("N83").Value = -27360
# Row 92
$ws.Range("H92").Value = 156.2
$ws.Range("I92").Value = 245.5
$ws.Range("J92").Value = 96.666664
$ws.Range("K92").Value = 736.5
$ws.Range("L92").Value = 289.999992
$ws.Range("M92").Value = 511.5
$ws.Range("N92").Value = -2785.999992
# Row 97
$ws.Range("H97").Value = 856
$ws.Range("J97").Value = 722.5
$ws.Range("L97").Value = 2167.5
$ws.Range("N97").Value = -3159.5
# Row 103
$ws.Range("H103").Value = 392.75
$ws.Range("I103").Value = 392.75
$ws.Range("K103").Value = 1178.25
$ws.Range("M103").Value = -299.25
# Row 107
$ws.Range("H107").Value = 894.125
$ws.Range("J107").Value = 755
$ws.Range("L107").Value = 2265
$ws.Range("N107").Value = -6105
# Row 109
$ws.Range("H109").Value = 3212.8572
$ws.Range("I109").Value = 1915
$ws.Range("K109").Value = 5745
$ws.Range("M109").Value = -4705
# Row 129
$ws.Range("H129").Value = 11004.667
$ws.Range("I129").Value = 1265
$ws.Range("J129").Value = 15874.5
$ws.Range("K129").Value = 3795
$ws.Range("L129").Value = 47623.5
$ws.Range("M129").Value = 1205
$ws.Range("N129").Value = -57623.5
# Row 131
$ws.Range("H131").Value = 1875.7142
$ws.Range("I131").Value = 1710
$ws.Range("K131").Value = 5130
$ws.Range("M131").Value = -90

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 20665
$ws.Range("J70").Value = 27497.5
$ws.Range("L70").Value = 27497.5
$ws.Range("N70").Value = -28037.5
# Row 73
$ws.Range("H73").Value = 20665
$ws.Range("J73").Value = 27497.5
$ws.Range("L73").Value = 27497.5
$ws.Range("N73").Value = -29369.5
# Row 102
$ws.Range("H102").Value = 850
$ws.Range("I102").Value = 700
$ws.Range("K102").Value = 700
$ws.Range("M102").Value = 922
# Row 107
$ws.Range("H107").Value = 999
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 999
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4839
# Row 122
$ws.Range("H122").Value = 1545.6
$ws.Range("I122").Value = 1545.6
$ws.Range("K122").Value = 4636.799999999999
$ws.Range("M122").Value = -2186.799999999999
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4999.5
$ws.Range("I16").Value = 4999.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4999.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4829.5
$ws.Range("N16").ClearContents()
# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
# Row 47
$ws.Range("H47").Value = 1200
$ws.Range("J47").Value = 1200
$ws.Range("L47").Value = 1200
$ws.Range("N47").Value = -2180
# Row 52
$ws.Range("H52").Value = 1200
$ws.Range("J52").Value = 1200
$ws.Range("L52").Value = 1200
$ws.Range("N52").Value = -1666
# Row 100
$ws.Range("H100").Value = 200
$ws.Range("I100").Value = 200
$ws.Range("K100").Value = 200
$ws.Range("M100").Value = 341

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 81
$ws.Range("H81").Value = 2500
$ws.Range("I81").Value = 2500
$ws.Range("K81").Value = 5000
$ws.Range("M81").Value = -3939
# Row 84
$ws.Range("H84").Value = 2500
$ws.Range("I84").Value = 2500
$ws.Range("K84").Value = 25000
$ws.Range("M84").Value = -19696
# Row 107
$ws.Range("H107").Value = 1742.1666
$ws.Range("J107").Value = 2225
$ws.Range("L107").Value = 6675
$ws.Range("N107").Value = -10515
# Row 122
$ws.Range("H122").Value = 1625
$ws.Range("I122").Value = 1625
$ws.Range("K122").Value = 4875
$ws.Range("M122").Value = -2425
# Row 136
$ws.Range("H136").Value = 2692.5
$ws.Range("I136").Value = 2692.5
$ws.Range("K136").Value = 8077.5
$ws.Range("M136").Value = -5527.5

Write-Output "Edit complete"